$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 1 (empty) and shift everything up so the data that was in
# rows 2-5 ends up in rows 1-4
$ws.Rows("1:1").Delete()

# Update the selection to match the target state
$ws.Range("F16").Select()
